$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at 491; this pushes the previous rows 491-515
# down to 492-516 unchanged, matching the target diff.
$ws.Rows.Item(491).Insert()

$ws.Cells.Item(491, 1).Value = 10
$ws.Cells.Item(491, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(491, 3).Value = "La Araucanía"
$ws.Cells.Item(491, 4).Value = 45147
$ws.Cells.Item(491, 5).Value = 9
$ws.Cells.Item(491, 6).Value = 100114013
$ws.Cells.Item(491, 7).Value = "Zanahoria"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 125
$ws.Cells.Item(491, 11).Value = 5000
$ws.Cells.Item(491, 12).Value = 5000
$ws.Cells.Item(491, 13).Value = 5000
$ws.Cells.Item(491, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(491, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(491, 16).Value = 200
$ws.Cells.Item(491, 17).Value = 25
$ws.Cells.Item(491, 18).Value = "Hortaliza"
